# Daily attendance processing - 2025-11-25 23:25:16
# Normalizes the "Recorded By" (column G) entries so that the "System"
# token is consistently ordered relative to the other recorder(s) listed
# in the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2

    if ($val -eq "system, System, backup@backdoor.com") {
        $cell.Value2 = "System, system, backup@backdoor.com"
    }
    elseif ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value2 = "dnasr281@gmail.com, System"
    }
}
